# Added New Mac-Address and Document Types
# Rebuild the "master-valid_document" lookup table on sheet1:
#  - Row 5 (PSP / POI) is replaced by DOC001 / POI
#  - 29 new document-type rows are appended below the existing 6 data rows
#  - the selection / used range grows from A1:F7 to A1:F36

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: PSP -> DOC001 (same doccat_code/lang_code/is_active/cr_by/cr_dtimes) ---
$ws.Range("A5").Value = "DOC001"

# --- New rows 8-36: doctyp_code, doccat_code, lang_code, is_active, cr_by, cr_dtimes ---
$newRows = @(
    @("DOC001", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC002", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC003", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC004", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC005", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC006", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC007", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC008", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC009", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC010", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC011", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC012", "POI", "ara", $true, "superadmin", "now()"),
    @("DOC001", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC013", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC014", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC015", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC004", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC005", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC006", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC016", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC017", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC018", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC008", "POA", "ara", $true, "superadmin", "now()"),
    @("DOC024", "POR", "ara", $true, "superadmin", "now()"),
    @("DOC025", "POR", "ara", $true, "superadmin", "now()"),
    @("DOC026", "POR", "ara", $true, "superadmin", "now()"),
    @("DOC001", "POR", "ara", $true, "superadmin", "now()"),
    @("DOC027", "POR", "ara", $true, "superadmin", "now()"),
    @("DOC028", "POR", "ara", $true, "superadmin", "now()")
)

$r = 8
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# --- Reset the "past the data" selection that Excel stores, now starting at G1 ---
$ws.Range("G1:XFD1048576").Select() | Out-Null
